# MHV-17222: bump the ValueSet metadata Version and Date stamps.
#
# The "Metadata" sheet (the workbook's active tab) holds a simple
# Property/Value table:
#   A3/B3 -> "Version" / "0.2.9-beta"          -> "0.2.10-beta"
#   A8/B8 -> "Date"    / "2023-02-16T09:21:54-06:00" -> "2023-12-06T12:46:33-06:00"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B3").Value = "0.2.10-beta"
$ws.Range("B8").Value = "2023-12-06T12:46:33-06:00"
